# Refresh the crypto-tracker symbol list: update Price (D), Volume(1h) (E)
# and Hora (G) columns for the rows whose source data changed, per the
# GitHub Actions run on Wed Jan 18 18:00:44 UTC 2023.
#
# Values in this sheet are stored as text (e.g. "293.59", "-2.73%", "18"),
# not numbers, so each cell is forced to a Text number format before the
# write (otherwise Excel would auto-convert numeric-looking strings like
# "293.59" or "18" into actual numbers/percentages) and the format is
# reset to the default "Normal" style afterwards so no visible formatting
# change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '293.59'
    'E2' = '-2.73%'
    'G2' = '18'
    'D3' = '31.33'
    'E3' = '-1.58%'
    'G3' = '18'
    'D4' = '4.968'
    'E4' = '-1.02%'
    'G4' = '18'
    'D5' = '0.07342'
    'E5' = '-6.16%'
    'G5' = '18'
    'D6' = '1.816'
    'E6' = '-10.74%'
    'G6' = '18'
    'D7' = '7.662'
    'E7' = '-1.53%'
    'G7' = '18'
    'D8' = '3.756'
    'E8' = '-0.39%'
    'G8' = '18'
    'D9' = '0.9090'
    'E9' = '-0.91%'
    'G9' = '18'
    'D10' = '0.1652'
    'E10' = '-4.98%'
    'G10' = '18'
    'D11' = '0.07585'
    'E11' = '-4.26%'
    'G11' = '18'
    'D12' = '0.08174'
    'E12' = '-6.40%'
    'G12' = '18'
    'D13' = '0.02986'
    'E13' = '-4.62%'
    'G13' = '18'
    'D14' = '0.09980'
    'E14' = '-0.23%'
    'G14' = '18'
    'D15' = '0.001494'
    'E15' = '-1.67%'
    'G15' = '18'
    'D16' = '0.005753'
    'E16' = '-1.14%'
    'G16' = '18'
    'G17' = '18'
    'E18' = '-0.13%'
    'G18' = '18'
    'D19' = '2.120'
    'E19' = '-6.51%'
    'G19' = '18'
    'D20' = '0.3277'
    'E20' = '0.13%'
    'G20' = '18'
    'D21' = '0.1307'
    'E21' = '1.31%'
    'G21' = '18'
    'D22' = '4.335'
    'E22' = '4.26%'
    'G22' = '18'
    'D23' = '0.1976'
    'E23' = '10.26%'
    'G23' = '18'
    'D24' = '0.04475'
    'E24' = '-2.53%'
    'G24' = '18'
    'D25' = '0.001228'
    'E25' = '-1.03%'
    'G25' = '18'
    'D26' = '0.004040'
    'E26' = '-9.70%'
    'G26' = '18'
    'D27' = '0.0001252'
    'E27' = '0.17%'
    'G27' = '18'
    'G28' = '18'
    'G29' = '18'
    'G30' = '18'
    'G31' = '18'
    'G32' = '18'
    'G33' = '18'
    'G34' = '18'
    'G35' = '18'
    'G36' = '18'
    'G37' = '18'
    'G38' = '18'
    'E39' = '-5.33%'
    'G39' = '18'
    'D40' = '0.04394'
    'E40' = '-7.44%'
    'G40' = '18'
    'D41' = '0.007442'
    'E41' = '0.53%'
    'G41' = '18'
    'D42' = '0.1320'
    'E42' = '-2.68%'
    'G42' = '18'
    'D43' = '0.002059'
    'E43' = '-3.82%'
    'G43' = '18'
    'D44' = '0.01108'
    'E44' = '3.18%'
    'G44' = '18'
    'D45' = '0.00005975'
    'E45' = '-1.36%'
    'G45' = '18'
    'D46' = '0.00000000751'
    'E46' = '0.16%'
    'G46' = '18'
    'D47' = '1.989'
    'E47' = '141.51%'
    'G47' = '18'
    'E48' = '-14.15%'
    'G48' = '18'
    'D49' = '0.00002104'
    'E49' = '0.16%'
    'G49' = '18'
    'D50' = '0.0002004'
    'E50' = '0.16%'
    'G50' = '18'
    'G51' = '18'
}

foreach ($ref in $updates.Keys) {
    $newValue = $updates[$ref]
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $newValue
    $cell.Style = 'Normal'
}

Write-Host "Updated $($updates.Count) cells"
